$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 34 questions and profiles updates: append new grad profile row.
# Fill order mirrors how the shared-string table ends up ordered upstream
# (username, then college, then full name, then the numeric counter).
$ws.Range("B18").Value = "sayalijadhav1101"
$ws.Range("D18").Value = "Pune University"
$ws.Range("A18").Value = "Sayali Mohan Jadhav"
$ws.Range("C18").Value = 0

# Match the formatting used on the prior (last) data row for the cells
# that inherit that row's style (name + counter columns); the username
# cell in that row keeps the workbook's default style, and the college
# cell keeps its own pre-existing style already set on D18.
$ws.Range("A17").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("C17").Copy()
$ws.Range("C18").PasteSpecial(-4122)

$ws.Range("D15").Select() | Out-Null
